# Replace the "YYYYQ4" quarter-label strings in column A (rows 2-39) with
# actual Excel date serials (Dec 31 of that year) and format them with a
# custom date/time number format, matching the source workbook's bugfix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerials = @(
    32142, 32508, 32873, 33238, 33603, 33969, 34334, 34699, 35064, 35430,
    35795, 36160, 36525, 36891, 37256, 37621, 37986, 38352, 38717, 39082,
    39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735,
    43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657
)

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dateSerials[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
